$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the status of the GCC/Accept Invalid/llvm_#30844 row from "Confirmed" to "Duplicate"
$ws.Range("E8").Value = "Duplicate"

# Move the active selection to F12 (was F15)
$ws.Range("F12").Select()
